$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.355.14"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "2.162.05"
$ws.Range("E3").Value = "  +3.37%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'228.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").Value = "'0.622"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("D7").Value = "'64.44"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.52%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'0.398"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.20%  "
$ws.Range("D10").Value = "'0.0859"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "'15.91"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("D13").Value = "2.485.04"
$ws.Range("E13").Value = "  +3.46%  "
$ws.Range("D14").Value = "'22.34"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.40%  "
$ws.Range("D15").Value = "'0.814"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'5.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.41%  "
$ws.Range("D17").Value = "2.179.23"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "39.305.93"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("D19").Value = "'72.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "'6.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "0.0₃0853"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "'231.45"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'2.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.66%  "
$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").Value = "'9.65"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.34%  "
$ws.Range("D27").Value = "'172.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "'0.138"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "'19.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.75%  "
$ws.Range("E30").Value = "  -2.76%  "
$ws.Range("D31").Value = "'2.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.74%  "
$ws.Range("E32").Value = "  +1.47%  "
$ws.Range("D33").Value = "'4.63"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.65%  "
$ws.Range("D34").Value = "'4.78"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").Value = "'7.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +8.61%  "
$ws.Range("D36").Value = "'0.0618"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("D37").Value = "'2.43"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("D38").Value = "'3.58"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.33%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'104.18"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.07%  "
$ws.Range("D41").Value = "'0.0229"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").Value = "'17.78"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "1.541.39"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("B44").Value = "FTXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D44").Value = "'4.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +5.67%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.18"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.95%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").Value = "'2.82"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("D47").Value = "'0.0924"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.93%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.10"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.17%  "
$ws.Range("D49").Value = "'7.85"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.50%  "
$ws.Range("D50").Value = "2.368.56"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").Value = "'2.97"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.02%  "
